$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new data rows after the existing item row (row 7), before the
# totals row (old row 8) / footer row (old row 9). ---
$ws.Rows("8:11").Insert()

# Copy the formatting (styles, borders, fonts, merges-pattern) of the first
# data row (row 7) down onto the four freshly-inserted rows.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q11").PasteSpecial(-4122)

# Match the row heights used by the new rows.
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5

# Re-create the per-row merges (A:B, C:G, H:K, L:M, N:O) for each new row,
# matching the pattern already used on row 7.
foreach ($r in 8..11) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

function Set-TextValue($range, $text) {
    # Force the value to be written as literal text (matching the shared
    # string cells used throughout this report), even though the cell's
    # number format looks numeric. Restore the original number format
    # afterwards so the cell's style index is unaffected.
    $fmt = $range.NumberFormat()
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $fmt
}

# Helper to fill one item row: item number, name, balance, order-limit,
# price, selling price, transaction count.
function Fill-ItemRow($r, $num, $name, $balance, $orderLimit, $price, $sellPrice, $txCount) {
    $ws.Range("A" + $r).Value = $num
    Set-TextValue $ws.Range("C" + $r) $name
    Set-TextValue $ws.Range("H" + $r) $balance
    Set-TextValue $ws.Range("L" + $r) $orderLimit
    Set-TextValue $ws.Range("N" + $r) $price
    Set-TextValue $ws.Range("P" + $r) $sellPrice
    Set-TextValue $ws.Range("Q" + $r) $txCount
}

Fill-ItemRow 8  2 "PANTAZOL 40MG 14 ENTERIC COATED TAB."       "0:1"   "1" "104.00"     "52.0000"    "0:1"
Fill-ItemRow 9  3 "TRIPLIXAM 10/2.5/10MG 15 F.C. TABS."        "-90:0" "1" "234.00"     "21294.0000" "91:0"
Fill-ItemRow 10 4 "VOLTAREN 75MG/3ML 3 AMP."                   "3:2"   "1" "51.00"      "16.8300"    "0:1"
Fill-ItemRow 11 5 "سرنجات 3 سم"                                 "0:0"   "0" "2.00"       "2.0000"     "1:0"

# The grand-total row (was row 8, now shifted to row 12) gets the new sum.
$ws.Range("P12").Value = 21393.830000000002

# The timestamp footer (was row 9, now shifted to row 13) gets the new
# generation time.
$ws.Range("A13").Value = "Wednesday, 3 September, 2025 9:56 AM"
